$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "29e Prix Agri-Center à Jettingen  "
$ws.Range("B8").Value = "27e Prix Super U Beaucourt  "
$ws.Range("B15").Value = "29e Grand Prix Gestimmo à Magstatt le bas  "
$ws.Range("B19").Value = "20e Grand Prix Gestimmo à Magstatt-le-Bas  "
$ws.Range("B30").Value = "12e Nuit des Gros Mollets à Flaxlanden  "
$ws.Range("B31").Value = "10e Nuit des Gros Mollets jeunes à Flaxlanden (poussins à minimes)  "
$ws.Range("B40").Value = "3e VTT Peugeot  "
$ws.Range("B41").Value = "24e Montée du Floridor, col du Hundsruck à Thann  "
$ws.Range("B45").Value = "11e Prix de Boron  "
$ws.Range("B46").Value = "15e Grimpée du Col Amic à Soultz"
$ws.Range("B50").Value = "5e VTT MS Automobile Rixheim  "
$ws.Range("B55").Value = "2e Cyclo-cross de Giromagny. Epreuve FFC ouverte aux FSGT"
$ws.Range("B65").Value = "2e cyclo-cross du Gloeckelsberg"
